# Business Case.pptx - "Added business case in .pdf"
#
# 1) Fix a missing space in the publisher-analysis commentary on slide 21
#    ("consistentlyand" -> "consistently and").
# 2) Add the small corner icon picture ("Imagen 6") to slide 21, matching
#    the identical icon picture already present (top-right corner) on
#    other slides of the deck (e.g. slide 1 / slide 6).

$p = $ppt.ActivePresentation

# --- 1. Fix the typo in the text box under the ribbon/funnel chart -------
# Only the missing space between "consistently" and "and" is touched; the
# rest of the run (including its trailing non-breaking space) is left
# completely untouched by doing a targeted substring replace instead of
# retyping the whole paragraph.
$slide21 = $p.Slides.Item(21)

for ($i = 1; $i -le $slide21.Shapes.Count; $i++) {
    $shp = $slide21.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*consistentlyand*") {
            $tr.Text = $tr.Text -replace "consistentlyand", "consistently and"
        }
    }
}

# --- 2. Add the corner icon picture to slide 21 ---------------------------
# The slide already has shape ids {1,2,4,5} in use, so the very first shape
# PowerPoint would hand out is id 3 (the only existing gap). The source
# deck's new picture is id 7, so two disposable placeholder shapes are
# created first (consuming ids 3 and 6) and removed again; the picture
# pasted afterwards then lands on id 7, matching the real authoring
# sequence that produced this deck.
$dummy1 = $slide21.Shapes.AddTextbox(1, 0, 0, 1, 1)
$dummy2 = $slide21.Shapes.AddTextbox(1, 0, 0, 1, 1)
$dummy1.Delete()
$dummy2.Delete()

# Re-use the exact same icon picture that is already pasted into the
# top-right corner of other slides in this deck (identical description,
# position and size), so the embedded image/relationship is consistent
# with the rest of the presentation.
$sourceSlide = $p.Slides.Item(1)
$sourceIcon = $null
for ($i = 1; $i -le $sourceSlide.Shapes.Count; $i++) {
    $cand = $sourceSlide.Shapes.Item($i)
    if ($cand.Type -eq 13) {
        $sourceIcon = $cand
    }
}

$sourceIcon.Copy()
$pastedRange = $slide21.Shapes.Paste()
$newPic = $pastedRange.Item(1)
$newPic.Name = "Imagen 6"

# NOTE: deliberately not touching Left/Top/Width/Height here - the pasted
# shape already inherits the exact same placement as the source icon
# (EMU 11379438, 65025, 745274, 733748), which is precisely where this
# corner icon belongs on every slide. Re-assigning those properties would
# round-trip the values through points and lose a little precision.
